# Einnahmen Old Oak Spez preise
# Replace the "Bob Marley: One Love" row (row 12) with a new "Fearless Flyers"
# entry, and add a missing "Minimal Abzug" value for "The Old Oak" (row 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verleiherabgaben")

# Row 12: new Suisa number, unchanged Minimal Abzug (150), new Abzug [%] of 30,
# new film title "Fearless Flyers", new Verleiher "Filmcoopi Zürich AG".
$ws.Range("B12").Value = "1018.463"
$ws.Range("D12").Value = 30
$ws.Range("F12").Value = "Fearless Flyers"
$ws.Range("G12").Value = "Filmcoopi Zürich AG"

# Row 15 ("The Old Oak"): add the missing Minimal Abzug value of 150.
$ws.Range("C15").Value = 150

# Leave the selection where the editor ended up after the edit.
$ws.Range("D16").Select() | Out-Null
